$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheet "Sheet 1" -> "Teachpoints"
# ---------------------------------------------------------------------------
$teach = $wb.Worksheets.Item(1)
$teach.Name = "Teachpoints"

# ---------------------------------------------------------------------------
# 2. Teachpoints restructuring: drop the old "Table 1" title row so that the
#    header row becomes row 1 and everything shifts up by one row.
# ---------------------------------------------------------------------------
$teach.Rows.Item(1).Delete() | Out-Null

# A few delay values moved around during the restructuring:
#  - left_pick_hover / left_pick / right_pick_hover / right_pick no longer
#    carry a value in column H (delay)
$teach.Range("H4").ClearContents()
$teach.Range("H5").ClearContents()
$teach.Range("H8").ClearContents()
$teach.Range("H9").ClearContents()

# - "rest" row gains a value of 700 in column F
$teach.Range("F13").Value = 700
$teach.Range("F13").NumberFormat = "@"

# Re-apply the freeze pane (header row + first column) now that the sheet
# starts at row 1 instead of row 2.
$teach.Activate()
$excel.ActiveWindow.FreezePanes = $false
$teach.Range("B2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Add the new "Sequences" sheet after Teachpoints, re-using the same
#    sheet-level formatting (page setup, gridlines, fonts, column widths...)
#    by duplicating Teachpoints and then replacing its contents.
# ---------------------------------------------------------------------------
$teach.Copy([System.Type]::Missing, $teach) | Out-Null
$seq = $wb.Worksheets.Item($teach.Index + 1)
$seq.Name = "Sequences"

# Clear all the copied data - we will rebuild it from scratch.
$seq.Cells.ClearContents() | Out-Null

# Drop the extra columns/rows copied from Teachpoints (sheet only needs
# columns A:E and rows 1:10).
$seq.Range("F1:H1").EntireColumn.Delete() | Out-Null
$seq.Range("A11:A19").EntireRow.Delete() | Out-Null

# Headers (row 1)
$seq.Range("A1").Value = "sequence"
$seq.Range("B1").Value = "teachpoints"
$seq.Range("C1").Value = "delays"
$seq.Range("D1").Value = "loop"

# Data (row 2)
$seq.Range("A2").Value = "Pick_Place_cups"
$seq.Range("B2").Value = "safety, left_pick_hover, left_pick, grip_closed, left_pick_hover, right_pick_hover, grip_open, right_pick"
$seq.Range("C2").Value = "1,1,1,1,1,1,1,1"
$seq.Range("D2").Value = "Yes"
$seq.Range("B2:D2").NumberFormat = "@"
$seq.Rows.Item(2).RowHeight = 104.55

# Re-apply the freeze pane for the new sheet as well.
$seq.Activate()
$excel.ActiveWindow.FreezePanes = $false
$seq.Range("B2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

$teach.Activate()
